$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = '4:00: Expo'
$ws.Range("D8").Value = '10:00: 1'
$ws.Range("E8").Value = '5:00: 1'
$ws.Range("F8").Value = '10:00: 1'
$ws.Range("I8").Value = '-'
$ws.Range("J8").Value = '-'
$ws.Range("M8").Value = '4:00: Expo'
$ws.Range("N8").Value = '10:00: 2'
$ws.Range("B9").Value = '10:00: 5'
$ws.Range("C9").Value = '4:00: 1'
$ws.Range("E9").Value = '4:00: 3'
$ws.Range("G9").Value = '5:00: 2'
$ws.Range("H9").Value = '10:00: 2'
$ws.Range("I9").Value = '5:00: 2'
$ws.Range("J9").Value = '10:00: 1'
$ws.Range("K9").Value = '5:00: 1'
$ws.Range("L9").Value = '10:00: 1'
$ws.Range("N9").Value = '10:00: 5'
$ws.Range("D10").Value = '-'
$ws.Range("H10").Value = '10:00: 4'
$ws.Range("I10").Value = '-'
$ws.Range("K10").Value = '4:00: 5'
$ws.Range("G11").Value = '4:00: 5'
$ws.Range("H11").Value = '-'
$ws.Range("L11").Value = '-'
$ws.Range("N11").Value = '-'
$ws.Range("O11").Value = '-'
$ws.Range("B12").Value = '-'
$ws.Range("E12").Value = '-'
$ws.Range("I12").Value = '-'
$ws.Range("J12").Value = '-'
$ws.Range("K12").Value = '-'
$ws.Range("L12").Value = '10:00: 5'
$ws.Range("B13").Value = '10:00: 1'
$ws.Range("F13").Value = '-'
$ws.Range("K13").Value = '4:00: 4'
$ws.Range("L13").Value = '-'
$ws.Range("M13").Value = '-'
$ws.Range("O13").Value = '4:00: 2'
$ws.Range("B14").Value = '-'
$ws.Range("H14").Value = '-'
$ws.Range("I14").Value = '5:00: 1'
$ws.Range("J14").Value = '-'
$ws.Range("K14").Value = '4:00: 3'
$ws.Range("M14").Value = '4:00: 1'
$ws.Range("O14").Value = '4:00: 1'
$ws.Range("F15").Value = '10:00: 2'
$ws.Range("J15").Value = '10:00: 3'
$ws.Range("K15").Value = '5:00: 2'
$ws.Range("L15").Value = '-'
$ws.Range("M15").Value = '4:00: 4'
$ws.Range("N15").Value = '10:00: 3'
$ws.Range("C16").Value = '-'
$ws.Range("E16").Value = '-'
$ws.Range("F16").Value = '10:00: 3'
$ws.Range("H16").Value = '10:00: 1'
$ws.Range("I16").Value = '-'
$ws.Range("M16").Value = '-'
$ws.Range("N16").Value = '10:00: 4'
$ws.Range("O16").Value = '4:00: 4'
$ws.Range("D17").Value = '10:00: 3'
$ws.Range("G17").Value = '-'
$ws.Range("H17").Value = '-'
$ws.Range("M17").Value = '4:00: 3'
$ws.Range("O17").Value = '4:00: 5'
$ws.Range("E18").Value = '4:00: 4'
$ws.Range("B19").Value = '-'
$ws.Range("E19").Value = '4:00: 5'
$ws.Range("I19").Value = '4:00: 4'
$ws.Range("J19").Value = '10:00: 4'
$ws.Range("O19").Value = '-'
$ws.Range("C20").Value = '4:00: 5'
$ws.Range("G20").Value = '4:00: 4'
$ws.Range("K20").Value = '-'
$ws.Range("L20").Value = '10:00: 2'
$ws.Range("B21").Value = '10:00: 4'
$ws.Range("C21").Value = '-'
$ws.Range("D21").Value = '10:00: 4'
$ws.Range("E21").Value = '-'
$ws.Range("F21").Value = '-'
$ws.Range("G21").Value = '-'
$ws.Range("M21").Value = '-'
$ws.Range("D22").Value = '-'
$ws.Range("E22").Value = '-'
$ws.Range("F22").Value = '10:00: 4'
$ws.Range("I22").Value = '4:00: 3'
$ws.Range("N22").Value = '10:00: 1'
$ws.Range("O22").Value = '4:00: 3'
$ws.Range("B23").Value = '10:00: 3'
$ws.Range("C23").Value = '4:00: 3'
$ws.Range("E23").Value = '-'
$ws.Range("G23").Value = '-'
$ws.Range("H23").Value = '-'
$ws.Range("L23").Value = '-'
$ws.Range("M23").Value = '4:00: 5'
$ws.Range("O23").Value = '-'
$ws.Range("C24").Value = '4:00: 2'
$ws.Range("E24").Value = '5:00: 2'
$ws.Range("F24").Value = '-'
$ws.Range("L24").Value = '10:00: 4'
$ws.Range("M24").Value = '4:00: 2'
$ws.Range("B25").Value = '-'
$ws.Range("D25").Value = '-'
$ws.Range("H25").Value = '10:00: 3'
$ws.Range("K25").Value = '-'
$ws.Range("L25").Value = '10:00: 3'
$ws.Range("M25").Value = '-'
$ws.Range("I26").Value = '-'
$ws.Range("M26").Value = '4:00: H/G'
$ws.Range("G27").Value = '4:00: H/G'
$ws.Range("I27").Value = '4:00: H/G'
$ws.Range("J27").Value = '11:00: H/G'
$ws.Range("M27").Value = '-'
$ws.Range("N27").Value = '11:00: H/G'
$ws.Range("L28").Value = '11:00: H/G'
$ws.Range("C29").Value = '-'
$ws.Range("H29").Value = '11:00: H/G'
$ws.Range("J29").Value = '-'
$ws.Range("L29").Value = '-'
$ws.Range("B30").Value = '11:00: H/G'
$ws.Range("C30").Value = '4:00: H/G'
$ws.Range("N30").Value = '-'
$ws.Range("B31").Value = '-'
$ws.Range("G31").Value = '-'
$ws.Range("H31").Value = '-'
$ws.Range("C32").Value = '-'
